$wb = $excel.ActiveWorkbook

# --- ProductLoanInput sheet ---
$wsInput = $wb.Worksheets.Item("ProductLoanInput")

# Row 1 (B1) - productname value: fix the product name string
$wsInput.Range("B1").Value = "773-RBI-EI-DB-DL-REC-NON-RNI-CTPD-SAR-MD-TR-1-LateRepayment"

# Row 18 (B18) - interestcalculationperiod value: Daily -> Same as repayment period
$wsInput.Range("B18").Value = "Same as repayment period"

# Move the active selection to B18 as in the edited workbook
$wsInput.Range("B18").Select()

# --- ProductLoanOutput sheet ---
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

# Row 1 (B1) - same corrected product name value
$wsOutput.Range("B1").Value = "773-RBI-EI-DB-DL-REC-NON-RNI-CTPD-SAR-MD-TR-1-LateRepayment"

# Select B1 on the output sheet to match the edited selection
$wsOutput.Range("B1").Select()

# Restore ProductLoanInput as the active/tab-selected sheet with its own
# selection on B18, matching the final state in the edited workbook.
$wsInput.Activate()
$wsInput.Range("B18").Select()
